$d = $word.ActiveDocument

$pairs = @(
    @('2025-09-16 Tuesday', '2025-09-17 Wednesday'),
    @('99-53=46', '85+12=97'),
    @('83-35=48', '36+10=46'),
    @('71+28=99', '45-20=25'),
    @('43+46=89', '89+2=91'),
    @('49+38=87', '63-0=63'),
    @('10+35=45', '55+7=62'),
    @('34+4=38', '99-41=58'),
    @('5+94=99', '80-3=77'),
    @('88-71=17', '10+81=91'),
    @('70-10=60', '94-91=3'),
    @('12+20=32', '45-11=34'),
    @('61+36=97', '17+28=45'),
    @('81-14=67', '42-9=33'),
    @('2+19=21', '34+50=84'),
    @('67-17=50', '33+39=72'),
    @('48+34=82', '82-19=63'),
    @('82-69=13', '81-68=13'),
    @('37-1=36', '51+10=61'),
    @('76-70=6', '16+82=98'),
    @('20+41=61', '71-67=4'),
    @('42+15=57', '2+26=28'),
    @('22+74=96', '72-51=21'),
    @('45+0=45', '93+5=98'),
    @('68-55=13', '13+27=40'),
    @('46-43=3', '96-46=50'),
    @('78+16=94', '15-15=0'),
    @('1+95=96', '36-19=17'),
    @('72+25=97', '41-35=6'),
    @('87-31=56', '70-34=36'),
    @('68+18=86', '33-0=33'),
    @('6+18=24', '30+36=66'),
    @('86-45=41', '56+10=66'),
    @('39+21=60', '5+44=49'),
    @('5-2=3', '49+44=93'),
    @('38+59=97', '57-13=44'),
    @('41-5=36', '66-17=49'),
    @('65-58=7', '9+30=39'),
    @('48-47=1', '29+60=89'),
    @('78-58=20', '63+2=65'),
    @('24+8=32', '2+84=86'),
    @('94-55=39', '19+22=41'),
    @('72-60=12', '37-35=2'),
    @('70-61=9', '80-74=6'),
    @('15+20=35', '90-21=69'),
    @('34+7=41', '51-22=29'),
    @('13+19=32', '76-2=74'),
    @('73+3=76', '97-23=74'),
    @('98-65=33', '79-63=16'),
    @('83-46=37', '22+47=69'),
    @('68-44=24', '98-81=17'),
    @('0+8=8', '28+18=46'),
    @('88-76=12', '56-18=38'),
    @('97-2=95', '74-48=26'),
    @('16+65=81', '52-13=39'),
    @('0+47=47', '32+8=40'),
    @('65+17=82', '92-0=92'),
    @('1+96=97', '21+54=75'),
    @('93-79=14', '51-41=10'),
    @('7+22=29', '13+60=73'),
    @('95-55=40', '59+29=88'),
    @('6+7=13', '49-45=4'),
    @('47-20=27', '96-95=1'),
    @('27+1=28', '65-21=44'),
    @('80-70=10', '38+25=63'),
    @('11+50=61', '65+23=88'),
    @('92+6=98', '73+5=78'),
    @('46-5=41', '55-0=55'),
    @('77-59=18', '64-18=46'),
    @('43+49=92', '73+21=94'),
    @('81-23=58', '41-1=40'),
    @('42+34=76', '40+12=52'),
    @('88-21=67', '51+44=95'),
    @('84-16=68', '55-54=1'),
    @('29+57=86', '26-16=10'),
    @('22+65=87', '41+47=88'),
    @('85-11=74', '95-1=94'),
    @('22+67=89', '30-11=19'),
    @('86+13=99', '14+77=91'),
    @('13+68=81', '3+8=11'),
    @('59+10=69', '99-82=17'),
    @('72+21=93', '25+16=41'),
    @('65+9=74', '57-10=47'),
    @('70-44=26', '61+6=67'),
    @('44-43=1', '72-13=59'),
    @('65+7=72', '98-35=63'),
    @('53-51=2', '52+7=59'),
    @('23-11=12', '98-92=6'),
    @('34+16=50', '27+19=46'),
    @('22-11=11', '67-4=63'),
    @('98-91=7', '26+19=45'),
    @('49+28=77', '28+55=83'),
    @('32-14=18', '52+46=98'),
    @('48+25=73', '61+15=76'),
    @('8+16=24', '11+26=37'),
    @('30+17=47', '6+53=59'),
    @('90-33=57', '50-33=17'),
    @('63+34=97', '43-10=33'),
    @('88-44=44', '89-25=64'),
    @('16+30=46', '33-26=7'),
    @('53-16=37', '96-83=13'),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()